$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.140.57"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.948.37"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "3.410.06"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "2.974.22"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +47.96%  "
$ws.Range("D19").Value = "51.081.84"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -6.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  +7.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("E31").Value = "  -4.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0442"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.272"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "1.988.99"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0324"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  +2.11%  "
